# Plasma_Gen_Issue Note.xlsx — issue-note format fix
#
# 1. Switch calculation mode from Manual back to Automatic.
# 2. On both the "Pipette" and "Transformer" sheets:
#    - delete the blank spacer row (row 2) so the table moves up one row
#    - bold the title cell in B1 (keeping its existing vertical-center alignment)
#    - refresh the zoom level / selected cell
#    - switch the print setup to landscape with a custom scale

$wb = $excel.ActiveWorkbook

# 1. Automatic calculation
$excel.Calculation = -4105   # xlCalculationAutomatic

# --- Pipette sheet -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("Pipette")
$ws1.Rows.Item(2).Delete()

$title1 = $ws1.Range("B1")
$title1.Font.Bold = $true
$title1.VerticalAlignment = -4108   # xlCenter

$ws1.Application.ActiveWindow.Zoom = 100
$ws1.Range("G17").Select()

$ws1.PageSetup.Orientation = 2      # xlLandscape
$ws1.PageSetup.Zoom = $false
$ws1.PageSetup.FitToPagesWide = $false
$ws1.PageSetup.FitToPagesTall = $false
$ws1.PageSetup.Zoom = 65

# --- Transformer sheet --------------------------------------------------
$ws2 = $wb.Worksheets.Item("Transformer")
$ws2.Rows.Item(2).Delete()

$title2 = $ws2.Range("B1")
$title2.Font.Bold = $true
$title2.VerticalAlignment = -4108   # xlCenter

$ws2.Application.ActiveWindow.Zoom = 100
$ws2.Range("H17").Select()

$ws2.PageSetup.Orientation = 2      # xlLandscape
$ws2.PageSetup.Zoom = $false
$ws2.PageSetup.FitToPagesWide = $false
$ws2.PageSetup.FitToPagesTall = $false
$ws2.PageSetup.Zoom = 73

$ws2.Activate()
